$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update player name and team for row 6 (5th ranked player)
$ws.Range("B6").Value = "Josh Hart"
$ws.Range("C6").Value = "NY"

# Update "Quantidade" (count) column D for rows 2-6
$ws.Range("D2").Value = 23
$ws.Range("D3").Value = 21
$ws.Range("D4").Value = 17
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 4
